$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row for SA5 assignment, matching style of existing rows (text formatted dates)
$ws.Range("A17").Value = "SA5"
$ws.Range("B17").Value = "October 06, 2024"

# Copy the number format / style used by the other deadline cells (e.g. B16) to B17
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Re-set the value after paste special formats, in case paste cleared it
$ws.Range("B17").Value = "October 06, 2024"

# Update selection to match post-edit state
$ws.Range("B18").Select() | Out-Null
